$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.408.43'
$ws.Range("E2").Value = '  +0.82%  '

$ws.Range("D3").Value = '3.377.58'
$ws.Range("E3").Value = '  +0.47%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '573.96'
$ws.Range("E5").Value = '  +0.68%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.37'
$ws.Range("E6").Value = '  +1.07%  '

$ws.Range("E7").Value = '  -0.07%  '

$ws.Range("D8").Value = '3.376.28'
$ws.Range("E8").Value = '  +0.48%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.474'
$ws.Range("E9").Value = '  -0.60%  '

$ws.Range("E10").Value = '  -1.06%  '

$ws.Range("E11").Value = '  +2.12%  '

$ws.Range("E12").Value = '  +0.04%  '

$ws.Range("D13").Value = '3.954.02'
$ws.Range("E13").Value = '  +0.46%  '

$ws.Range("E14").Value = '  +2.44%  '

$ws.Range("E15").Value = '  +2.68%  '

$ws.Range("B16").Value = 'Avalanche'
$ws.Range("C16").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.88'
$ws.Range("E16").Value = '  +3.48%  '

$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.380.84'
$ws.Range("E17").Value = '  +0.62%  '

$ws.Range("D18").Value = '61.469.33'
$ws.Range("E18").Value = '  +0.76%  '

$ws.Range("E19").Value = '  +0.75%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.85'
$ws.Range("E20").Value = '  +1.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.38'
$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '376.31'
$ws.Range("E22").Value = '  +0.56%  '

$ws.Range("E23").Value = '  -2.98%  '

$ws.Range("D24").Value = '3.522.49'
$ws.Range("E24").Value = '  +0.76%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.00'
$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("E26").Value = '  +7.45%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '71.24'
$ws.Range("E27").Value = '  +0.89%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.73'
$ws.Range("E28").Value = '  +4.79%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.48'
$ws.Range("E29").Value = '  -2.75%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("E30").Value = '  -0.01%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.161'
$ws.Range("E31").Value = '  +4.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.25'
$ws.Range("E32").Value = '  +2.10%  '

$ws.Range("E33").Value = '  +1.85%  '

$ws.Range("E34").Value = '  +0.03%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.47'
$ws.Range("E35").Value = '  +0.16%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.27'
$ws.Range("E36").Value = '  -4.84%  '

$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.82'
$ws.Range("E37").Value = '  -0.96%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.54'
$ws.Range("E38").Value = '  -0.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '165.43'
$ws.Range("E39").Value = '  +1.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0775'
$ws.Range("E40").Value = '  -1.51%  '

$ws.Range("E41").Value = '  +0.07%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.774'
$ws.Range("E42").Value = '  +2.40%  '

$ws.Range("B43").Value = 'Stacks'
$ws.Range("C43").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.72'
$ws.Range("E43").Value = '  +7.47%  '

$ws.Range("B44").Value = 'ONDO'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.22'
$ws.Range("E44").Value = '  +1.34%  '

$ws.Range("E45").Value = '  +0.54%  '

$ws.Range("E46").Value = '  +0.11%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.52'
$ws.Range("E47").Value = '  +8.04%  '

$ws.Range("E48").Value = '  -1.84%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.71'
$ws.Range("E49").Value = '  -1.55%  '

$ws.Range("D50").Value = '2.348.01'
$ws.Range("E50").Value = '  +4.98%  '

$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0261'
$ws.Range("E51").Value = '  +1.57%  '
